$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 44.29505033333334
$ws.Range("H2").Value = 132.885151
$ws.Range("I2").Value = 0.9830698162761968
$ws.Range("J2").Value = 0.9830698162761969
$ws.Range("M2").Value = 1.910418
$ws.Range("N2").Value = 5.731254
$ws.Range("O2").Value = 0.01809124304049503
$ws.Range("P2").Value = 0.01809124304049503
$ws.Range("Q2").Value = 84.622061467706
$ws.Range("R2").Value = 761.598553209354
$ws.Range("S2").Value = 0.01778495497202747
$ws.Range("T2").Value = 0.01778495497202747
$ws.Range("G3").Value = 44.29505033333334
$ws.Range("H3").Value = 132.885151
$ws.Range("I3").Value = 0.9830698162761968
$ws.Range("J3").Value = 0.9830698162761969
$ws.Range("O3").Value = 0.302988173785169
$ws.Range("P3").Value = 0.302988173785169
$ws.Range("Q3").Value = 1417.231740718187
$ws.Range("R3").Value = 12755.08566646369
$ws.Range("S3").Value = 0.2978585283368465
$ws.Range("T3").Value = 0.2978585283368465
$ws.Range("G4").Value = 44.29505033333334
$ws.Range("H4").Value = 132.885151
$ws.Range("I4").Value = 0.9830698162761968
$ws.Range("J4").Value = 0.9830698162761969
$ws.Range("M4").Value = 37.858701
$ws.Range("N4").Value = 113.576103
$ws.Range("O4").Value = 0.3585136661130873
$ws.Range("P4").Value = 0.3585136661130873
$ws.Range("Q4").Value = 1676.953066349617
$ws.Range("R4").Value = 15092.57759714655
$ws.Range("S4").Value = 0.3524439638782985
$ws.Range("T4").Value = 0.3524439638782986
$ws.Range("G5").Value = 44.29505033333334
$ws.Range("H5").Value = 132.885151
$ws.Range("I5").Value = 0.9830698162761968
$ws.Range("J5").Value = 0.9830698162761969
$ws.Range("M5").Value = 33.83466466666667
$ws.Range("N5").Value = 101.503994
$ws.Range("O5").Value = 0.3204069170612486
$ws.Range("P5").Value = 0.3204069170612486
$ws.Range("Q5").Value = 1498.708174421455
$ws.Range("R5").Value = 13488.3735697931
$ws.Range("S5").Value = 0.3149823690890243
$ws.Range("T5").Value = 0.3149823690890243
$ws.Range("I6").Value = 0.006814145293655052
$ws.Range("J6").Value = 0.006814145293655053
$ws.Range("M6").Value = 1.910418
$ws.Range("N6").Value = 5.731254
$ws.Range("O6").Value = 0.01809124304049503
$ws.Range("P6").Value = 0.01809124304049503
$ws.Range("Q6").Value = 0.5865575489579999
$ws.Range("R6").Value = 5.279017940622
$ws.Range("S6").Value = 0.0001232763586207589
$ws.Range("T6").Value = 0.0001232763586207589
$ws.Range("I7").Value = 0.006814145293655052
$ws.Range("J7").Value = 0.006814145293655053
$ws.Range("O7").Value = 0.302988173785169
$ws.Range("P7").Value = 0.302988173785169
$ws.Range("S7").Value = 0.002064605438431349
$ws.Range("T7").Value = 0.002064605438431349
$ws.Range("I8").Value = 0.006814145293655052
$ws.Range("J8").Value = 0.006814145293655053
$ws.Range("M8").Value = 37.858701
$ws.Range("N8").Value = 113.576103
$ws.Range("O8").Value = 0.3585136661130873
$ws.Range("P8").Value = 0.3585136661130873
$ws.Range("Q8").Value = 11.623794826731
$ws.Range("R8").Value = 104.614153440579
$ws.Range("S8").Value = 0.002442964210655513
$ws.Range("T8").Value = 0.002442964210655513
$ws.Range("I9").Value = 0.006814145293655052
$ws.Range("J9").Value = 0.006814145293655053
$ws.Range("M9").Value = 33.83466466666667
$ws.Range("N9").Value = 101.503994
$ws.Range("O9").Value = 0.3204069170612486
$ws.Range("P9").Value = 0.3204069170612486
$ws.Range("Q9").Value = 10.38829092727133
$ws.Range("R9").Value = 93.49461834544201
$ws.Range("S9").Value = 0.002183299285947432
$ws.Range("T9").Value = 0.002183299285947432
$ws.Range("G10").Value = 0.3685326666666667
$ws.Range("H10").Value = 1.105598
$ws.Range("I10").Value = 0.008179093108268589
$ws.Range("J10").Value = 0.008179093108268589
$ws.Range("M10").Value = 1.910418
$ws.Range("N10").Value = 5.731254
$ws.Range("O10").Value = 0.01809124304049503
$ws.Range("P10").Value = 0.01809124304049503
$ws.Range("Q10").Value = 0.704051439988
$ws.Range("R10").Value = 6.336462959892001
$ws.Range("S10").Value = 0.000147969961272525
$ws.Range("T10").Value = 0.000147969961272525
$ws.Range("G11").Value = 0.3685326666666667
$ws.Range("H11").Value = 1.105598
$ws.Range("I11").Value = 0.008179093108268589
$ws.Range("J11").Value = 0.008179093108268589
$ws.Range("O11").Value = 0.302988173785169
$ws.Range("P11").Value = 0.302988173785169
$ws.Range("Q11").Value = 11.79129922555867
$ws.Range("R11").Value = 106.121693030028
$ws.Range("S11").Value = 0.002478168484093161
$ws.Range("T11").Value = 0.002478168484093161
$ws.Range("G12").Value = 0.3685326666666667
$ws.Range("H12").Value = 1.105598
$ws.Range("I12").Value = 0.008179093108268589
$ws.Range("J12").Value = 0.008179093108268589
$ws.Range("M12").Value = 37.858701
$ws.Range("N12").Value = 113.576103
$ws.Range("O12").Value = 0.3585136661130873
$ws.Range("P12").Value = 0.3585136661130873
$ws.Range("Q12").Value = 13.952168036066
$ws.Range("R12").Value = 125.569512324594
$ws.Range("S12").Value = 0.002932316655725659
$ws.Range("T12").Value = 0.002932316655725659
$ws.Range("G13").Value = 0.3685326666666667
$ws.Range("H13").Value = 1.105598
$ws.Range("I13").Value = 0.008179093108268589
$ws.Range("J13").Value = 0.008179093108268589
$ws.Range("M13").Value = 33.83466466666667
$ws.Range("N13").Value = 101.503994
$ws.Range("O13").Value = 0.3204069170612486
$ws.Range("P13").Value = 0.3204069170612486
$ws.Range("Q13").Value = 12.46917919537911
$ws.Range("R13").Value = 112.222612758412
$ws.Range("S13").Value = 0.002620638007177243
$ws.Range("T13").Value = 0.002620638007177243
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.08727466666666667
$ws.Range("H14").Value = 0.261824
$ws.Range("I14").Value = 0.001936945321879485
$ws.Range("J14").Value = 0.001936945321879485
$ws.Range("M14").Value = 1.910418
$ws.Range("N14").Value = 5.731254
$ws.Range("O14").Value = 0.01809124304049503
$ws.Range("P14").Value = 0.01809124304049503
$ws.Range("Q14").Value = 0.166731094144
$ws.Range("R14").Value = 1.500579847296
$ws.Range("S14").Value = 0.00003504174857427163
$ws.Range("T14").Value = 0.00003504174857427163
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.08727466666666667
$ws.Range("H15").Value = 0.261824
$ws.Range("I15").Value = 0.001936945321879485
$ws.Range("J15").Value = 0.001936945321879485
$ws.Range("O15").Value = 0.302988173785169
$ws.Range("P15").Value = 0.302988173785169
$ws.Range("Q15").Value = 2.792375825962667
$ws.Range("R15").Value = 25.131382433664
$ws.Range("S15").Value = 0.0005868715257979915
$ws.Range("T15").Value = 0.0005868715257979916
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.08727466666666667
$ws.Range("H16").Value = 0.261824
$ws.Range("I16").Value = 0.001936945321879485
$ws.Range("J16").Value = 0.001936945321879485
$ws.Range("M16").Value = 37.858701
$ws.Range("N16").Value = 113.576103
$ws.Range("O16").Value = 0.3585136661130873
$ws.Range("P16").Value = 0.3585136661130873
$ws.Range("Q16").Value = 3.304105510208
$ws.Range("R16").Value = 29.736949591872
$ws.Range("S16").Value = 0.0006944213684076081
$ws.Range("T16").Value = 0.0006944213684076082
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.08727466666666667
$ws.Range("H17").Value = 0.261824
$ws.Range("I17").Value = 0.001936945321879485
$ws.Range("J17").Value = 0.001936945321879485
$ws.Range("M17").Value = 33.83466466666667
$ws.Range("N17").Value = 101.503994
$ws.Range("O17").Value = 0.3204069170612486
$ws.Range("P17").Value = 0.3204069170612486
$ws.Range("Q17").Value = 2.952909080561778
$ws.Range("R17").Value = 26.576181725056
$ws.Range("S17").Value = 0.0006206106790996135
$ws.Range("T17").Value = 0.0006206106790996136

Write-Host "Applied 182 cell updates"
